$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddStudent")

# Fill in values in the order the shared strings are introduced:
# Student1, Test1, Test2, Student2
$ws.Range("A2").Value = "Student1"
$ws.Range("C2").Value = "Test1"
$ws.Range("D2").Value = "Test1"
$ws.Range("E2").Value = "Test1"
$ws.Range("C3").Value = "Test2"
$ws.Range("D3").Value = "Test2"
$ws.Range("E3").Value = "Test2"
$ws.Range("A3").Value = "Student2"

$ws.Range("B2").Value = 123456
$ws.Range("F2").Value = 123456
$ws.Range("B3").Value = 123456
$ws.Range("F3").Value = 123456

# Widen column D slightly (target stored width 16.28515625; the COM layer
# quantizes ColumnWidth to 1/6-character steps, so 15.5 is the input that
# lands closest to the target on save)
$ws.Columns.Item(4).ColumnWidth = 15.5

# Move the active selection to C3
$ws.Range("C3").Select()
